$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of column G (Recorded By) for data rows, but keep formatting
$ws.Range("G2:G115").ClearContents()

# Adjust column G width from 50 to 13
# (ColumnWidth undergoes pixel-grid quantization on save; 12.15 lands exactly on width=13)
$ws.Columns.Item(7).ColumnWidth = 12.15
